$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("AI3").Value = 285
$ws3.Range("AJ3").Value = 280
$ws3.Range("AK3").Value = 293
$ws3.Range("AL3").Value = 292
$ws3.Range("AN3").Value = 285
$ws3.Range("AI4").Value = 777
$ws3.Range("AJ4").Value = 817
$ws3.Range("AK4").Value = 740
$ws3.Range("AL4").Value = 754
$ws3.Range("AM4").Value = 784
$ws3.Range("AN4").Value = 794
$ws3.Range("AI6").Value = 5
$ws3.Range("AJ6").Value = 26
$ws3.Range("AK6").Value = 32
$ws3.Range("AL6").Value = 31
$ws3.Range("AM6").Value = 46
$ws3.Range("AN6").Value = 10
$ws3.Range("AI7").Value = 159
$ws3.Range("AJ7").Value = 86
$ws3.Range("AK7").Value = 146
$ws3.Range("AL7").Value = 121
$ws3.Range("AM7").Value = 77
$ws3.Range("AN7").Value = 111
$ws3.Range("AI8").Value = 49
$ws3.Range("AJ8").Value = 75
$ws3.Range("AK8").Value = 56
$ws3.Range("AL8").Value = 92
$ws3.Range("AM8").Value = 91
$ws3.Range("AN8").Value = 109
$ws3.Range("AJ9").Value = 71
$ws3.Range("AK9").Value = 47
$ws3.Range("AL9").Value = 35
$ws3.Range("AM9").Value = 54
$ws3.Range("AN9").Value = 40
$ws3.Range("AI10").Value = 19
$ws3.Range("AJ10").Value = 22
$ws3.Range("AK10").Value = 12
$ws3.Range("AL10").Value = 13
$ws3.Range("AM10").Value = 19
$ws3.Range("J14").Value = 287
$ws3.Range("J15").Value = 778
$ws3.Range("J17").Value = 117
$ws3.Range("J18").Value = 79
$ws3.Range("J23").Value = 294
$ws3.Range("J24").Value = 690
$ws3.Range("J25").Value = 16
$ws3.Range("J32").Value = 582
$ws3.Range("J33").Value = 1970
$ws3.Range("J34").Value = 48

$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("J19").Value = 406
$ws4.Range("J20").Value = 344
$ws4.Range("J21").Value = 326
$ws4.Range("J22").Value = 272
$ws4.Range("J23").Value = 273
$ws4.Range("J24").Value = 101
$ws4.Range("J28").Value = 1095
$ws4.Range("J29").Value = 946
$ws4.Range("J30").Value = 838
$ws4.Range("J31").Value = 753
$ws4.Range("J32").Value = 755
$ws4.Range("J33").Value = 279

$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("G2").Value = 26
$ws5.Range("G3").Value = 14
$ws5.Range("G4").Value = 48
$ws5.Range("G5").Value = 151
$ws5.Range("G6").Value = 199
$ws5.Range("G7").Value = 117
$ws5.Range("G8").Value = 17
$ws5.Range("G9").Value = 27
$ws5.Range("G10").Value = 34
$ws5.Range("G11").Value = 2
$ws5.Range("G12").Value = 37
$ws5.Range("G13").Value = 3
$ws5.Range("G14").Value = 18
$ws5.Range("G16").Value = 0
$ws5.Range("G17").Value = 53
$ws5.Range("G19").Value = 13
$ws5.Range("G21").Value = 216
$ws5.Range("G22").Value = 25
$ws5.Range("G24").Value = 0
$ws5.Range("G27").Value = 0
$ws5.Range("G28").Value = 0

$ws7 = $wb.Worksheets.Item("Sheet7")
$ws7.Range("G2").Value = 8.800000000000001
$ws7.Range("G5").Value = 59.74
